$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2017-02-15 06:06:13"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2017-02-15 06:05:54"
$wsZhCn.Range("L2").Value = "2017-02-15 06:06:40"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2017-02-15 06:06:13"
$wsDeDe.Range("L2").Value = "2017-02-15 06:07:07"
